# Applies the price/volume/coin-name updates described by the commit diff.
# Numeric-looking strings (prices, percentages) are written while the cell is
# temporarily formatted as Text ("@") so Excel keeps them as literal strings
# (matching the original inlineStr cell typing) instead of coercing them to
# numbers/dates/percentages. The style is then reset to "Normal" so no stray
# cell formatting is left behind, keeping cells visually identical to before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $range = $ws.Range($addr)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue "D2" "327.95"
Set-TextValue "E2" "-0.58%"
Set-TextValue "D3" "39.52"
Set-TextValue "E3" "-1.63%"
Set-TextValue "D4" "5.710"
Set-TextValue "E4" "5.74%"
Set-TextValue "D5" "0.08049"
Set-TextValue "E5" "-1.02%"
Set-TextValue "D6" "2.012"
Set-TextValue "E6" "4.70%"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue "D7" "4.495"
Set-TextValue "E7" "-0.69%"
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
Set-TextValue "D8" "8.633"
Set-TextValue "E8" "-0.25%"
Set-TextValue "D10" "0.9235"
Set-TextValue "E10" "-2.35%"
Set-TextValue "D11" "0.1261"
Set-TextValue "E11" "-7.23%"
Set-TextValue "D12" "0.1963"
Set-TextValue "E12" "-0.89%"
Set-TextValue "D13" "8.754"
Set-TextValue "E13" "21.84%"
Set-TextValue "D14" "0.09185"
Set-TextValue "E14" "-1.55%"
Set-TextValue "D15" "0.03565"
Set-TextValue "E15" "0.33%"
Set-TextValue "D16" "0.1051"
Set-TextValue "E16" "9.54%"
Set-TextValue "D17" "0.001302"
Set-TextValue "E17" "-1.42%"
Set-TextValue "D18" "0.006374"
Set-TextValue "E18" "-3.78%"
Set-TextValue "D19" "3.365"
Set-TextValue "E20" "-1.27%"
Set-TextValue "D21" "0.1350"
Set-TextValue "E21" "1.44%"
Set-TextValue "E22" "10.14%"
Set-TextValue "D23" "0.04401"
Set-TextValue "E23" "-0.64%"
Set-TextValue "E24" "3.05%"
Set-TextValue "E25" "7.21%"
Set-TextValue "D26" "0.0001187"
Set-TextValue "E26" "-1.10%"
Set-TextValue "D39" "0.02495"
Set-TextValue "E39" "-0.31%"
Set-TextValue "D40" "0.05349"
Set-TextValue "E40" "2.37%"
Set-TextValue "D41" "0.007485"
Set-TextValue "E41" "-1.52%"
Set-TextValue "D42" "0.009920"
Set-TextValue "E42" "9.05%"
Set-TextValue "D43" "0.1408"
Set-TextValue "E43" "-1.46%"
Set-TextValue "D44" "0.002111"
Set-TextValue "E44" "-2.35%"
Set-TextValue "D45" "0.01008"
Set-TextValue "E45" "-6.33%"
Set-TextValue "D46" "0.00006684"
Set-TextValue "E46" "1.43%"
Set-TextValue "E47" "-0.26%"
Set-TextValue "D48" "0.003035"
Set-TextValue "E48" "-9.27%"
Set-TextValue "D49" "0.002278"
Set-TextValue "E49" "-5.14%"
Set-TextValue "D50" "0.00002096"
Set-TextValue "E50" "-0.26%"
Set-TextValue "D51" "0.0001996"
Set-TextValue "E51" "-0.26%"
